# Release implementations for main, init_lister and get_source_line functions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phase 2 Tasks")

# Row 2: main() -> now released, % completed 80% -> 90%, test driver required "Yes" -> "Yes (0%)"
$ws.Range("F2").Value = 0.9
$ws.Range("G2").Value = "Yes (0%)"

# Row 4: get_source_line() -> now released, % completed 55% -> 100%, test driver required "Yes" -> "Yes (0%)"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "Yes (0%)"

# Column G widened (manually sized, no longer auto "best fit")
$ws.Columns("G").ColumnWidth = 15.86

# Remove now-empty trailing rows 11-14
$ws.Rows("11:14").Delete()

# Move the active selection
[void]$ws.Range("D4").Select()
